# IFCB_beadvoltages_temp.xlsx — add new laser-alignment signal readings to IFCB10
# and tidy up the active-sheet/active-tab bookkeeping.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("IFCB10")
$ws2 = $wb.Worksheets.Item("IFCB1")

# ---------------------------------------------------------------------------
# 1. IFCB10 (sheet1): insert the new "hvB" column (C) and the two new
#    "PMTA calc peak" / "PMTA calc int" columns (which land at I:J once C
#    has been inserted), pushing the old C:M block right.
# ---------------------------------------------------------------------------
$ws1.Columns("C").Insert()
$ws1.Columns("I:J").Insert()

# Column header text (C1 reuses the existing "hvB" shared string).
$ws1.Range("C1").Value = "hvB"

# I1/J1 need their new header text entered in this order so the new shared
# strings land at the indices the workbook expects.
$ws1.Range("I1").Value = "PMTA`ncalc peak"
$ws1.Range("J1").Value = "PMTA`ncalc int"

# ---------------------------------------------------------------------------
# 2. Append the two new data rows (14 & 15) captured on 2016-03-21 / 03-23.
# ---------------------------------------------------------------------------

# Row 14
$ws1.Range("A14").Value = 20160321
$ws1.Range("B14").Value = 201207
$ws1.Range("C14").Value = 0.7
$ws1.Range("D14").Value = "NA"
$ws1.Range("E14").Value = 1.3754
$ws1.Range("F14").Value = 0.1227
$ws1.Range("H14").Value = "1.2-1.6"
$ws1.Range("G14").Value = "0.1-0.2"
$ws1.Range("I14").Value = 3.4763000000000002
$ws1.Range("J14").Value = 0.30814000000000002
$ws1.Range("K14").Value = "NAN"
$ws1.Range("L14").Value = "NAN"
$ws1.Range("M14").Value = "NAN"
$ws1.Range("N14").Value = 0.11645
$ws1.Range("O14").Value = "completely realigned laser including move PMTS because camera at edge of adjustment and delay at 10us, move laser down when looking at telescope"
$ws1.Range("P14").Value = "pump1 us all signals"

# Row 15
$ws1.Range("A15").Value = 20160323
$ws1.Range("B15").Value = 212722
$ws1.Range("C15").Value = 0.7
$ws1.Range("D15").Value = "NA"
$ws1.Range("E15").Value = 0.93152999999999997
$ws1.Range("F15").Value = 0.10829
$ws1.Range("G15").Value = "0.07-0.13"
$ws1.Range("H15").Value = "0.8-1.15"
$ws1.Range("I15").Value = 2.9788999999999999
$ws1.Range("J15").Value = 0.28516000000000002
$ws1.Range("K15").Value = 0.20415
$ws1.Range("L15").Value = 0.10391
$ws1.Range("M15").Value = 0.25572
$ws1.Range("N15").Value = 0.086920999999999998
$ws1.Range("O15").Value = "selected points single beads"
$ws1.Range("P15").Value = "adjusted laser vert,horz,foc again after moving needle and then camera stack"

# ---------------------------------------------------------------------------
# 3. View/selection bookkeeping: IFCB10 becomes the active tab, gets a
#    frozen header row, and the cursor ends up on the newly added A14; the
#    IFCB1 sheet loses its old tabSelected flag and its own selection moves.
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("A14").Select()

$ws2.Activate()
$ws2.Range("N17").Select()

$ws1.Activate()
